$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Sheet 1")
$wsRefs = $wb.Worksheets.Item("References")

$wsData.Range("E5").Value = "Based on official disease reports to the WOAH"
$wsData.Range("E6").Value = "Glanders is a disease listed in the World Organisation for Animal Health ({ref009:WOAH}) Terrestrial Animal Health Code and must be reported to the WOAH. The map to the right displays outbreak points reported to the WOAH early warning system since 2005."
$wsData.Range("E7").Value = "As described in the WOAH {ref005:Terrestrial Animal Health Code}, the WOAH early warning system includes immediate notifications and follow-up reports on:"
$wsData.Range("E14").Value = "Countries are coloured according to the available information regarding their stable disease situation (disease status legend). This information is provided by countries through the WOAH monitoring system, which is a different reporting channel.<br>Immediate notifications (points) and disease status (country/region colours) are reported to the WOAH in different spatial and temporal scales, and therefore are displayed in the map as layers which can be filtered independently."
$wsData.Range("E17").Value = "For more up to date reports, visit the original data source: {ref001:WOAH-WAHIS}."
$wsData.Range("E21").Value = "Glanders is a disease caused by the bacteria <i>Burkholderia mallei</i>. Many animal species are known to be susceptible but disease generally occur in horses, donkeys and mules. Infections may cause ulcerations and nodules on the skin and in the respiratory tract and there are four different clinical presentations. The nasal and pulmonay forms, where clincical signs include nasal discharge, cough, fever and dyspnoea. These forms tend to be acute. A cutaneous form, also referred to as farcy, also exists where infected animals develop enlarged lymph nodes, nodular abscesses or dry ulcers. This form tend to be more chronic in nature. Nodules may also form in other internal organs, subsequently leading to wasting and eventually death.  Asymptomatic carriers are also possible.  Donkeys and mules tend to develop the acute respiratory forms and horses more often develop the chronic or asymtomatic forms ({ref008:WOAH})."
$wsData.Range("E42").Value = "<i>Burkholderia mallei</i> is a Gram-negative aerobic bacterium {ref008:WOAH Technical disease card} "
$wsData.Range("E54").Value = "Transmission occurs through contaminated feed or water, through aerosols or via fomites such as riding tack. Dermal transmission through abarasions in the skin or through the mucosal membranes is also possible. Over-crowding and poor sanitation predisposes are known risk factors for disease transmission. The incubation period varies from a few days to several months {ref008:WOAH Technical disease card}."
$wsData.Range("E66").Value = "Glanders can be diagnosed using by identifying the disease causing bacteria using either bacterial culture or Polymerase Chain Reaction (PCR). Serological tests such as the Complement fixation test (CFT) is also available and infected animals test positive on CFT approximately one week following infection. Enzyme-linked immunosorbent assays (ELISA) and Immunoblot assays are available but are not yet fully validated, ({ref009:WOAH, Terrestrial Manual})."
$wsData.Range("E84").Value = "For further information, visit the {ref034:CDC page on Glanders}), or the WOAH (panel to the right)."
$wsData.Range("E95").Value = "Geographical distribution data has been kindly provided by the World Organisation of Animal Health (WOAH). {ref001:WOAH-WAHIS} (WOAH World Animal Health Information System) is the original source of these data."

$wsRefs.Range("C2").Value = "WOAH-WAHIS (WOAH World Animal Health Information System)"
$wsRefs.Range("C6").Value = "WOAH (World Organisation for Animal Health). Terrestrial Animal Health Code 2021. WOAH, Paris, France"
$wsRefs.Range("C9").Value = "WOAH (World Organisation for Animal Health) Technical Disease Card: Glanders. 2021."
$wsRefs.Range("C10").Value = "WOAH (World Organisation for Animal Health), 2021. Glanders. Chapter 3.6.11. WOAH Terrestrial Animal Health Code 2021. WOAH, Paris, France"
$wsRefs.Range("C11").Value = "WOAH (World Organisation for Animal Health), 2021, Glanders, WOAH, Paris, France"
